$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two oldest quarters (2020Q4 and 2021Q1); this shifts all the
# following rows up by two and drops their now-unused shared strings.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# A handful of the remaining rows were recomputed with updated source data
# as part of finishing the script refactor; apply the corrected values.
$ws.Range("C3").Value = 206
$ws.Range("E3").Value = 110

$ws.Range("C4").Value = 214
$ws.Range("D4").Value = 124
$ws.Range("E4").Value = 90
$ws.Range("F4").Value = 60.19417475728155

$ws.Range("F5").Value = 60.74766355140186
